# Fix typos in main block assignments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct mistaken "bead_g" entries that should have been "bead_b"
$ws.Range("B14").Value = "stimuli/bead_b.PNG"
$ws.Range("C15").Value = "stimuli/bead_b.PNG"
$ws.Range("D16").Value = "stimuli/bead_b.PNG"

# Update the selection / scroll position left in the sheet view
$ws.Range("D16").Select()
